$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2086752637749121
$ws.Range("C2").Value = 0.52989449003517
$ws.Range("J2").Value = 0.009378663540445486
$ws.Range("P2").Value = 0.1371629542790152
$ws.Range("S2").Value = 0.1148886283704572
$ws.Range("B3").Value = 0.00211864406779661
$ws.Range("C3").Value = 0.02754237288135593
$ws.Range("J3").Value = 0.05296610169491525
$ws.Range("P3").Value = 0.7076271186440678
$ws.Range("S3").Value = 0.2097457627118644
$ws.Range("J4").Value = 0.05343511450381679
$ws.Range("P4").Value = 0.6870229007633588
$ws.Range("S4").Value = 0.2595419847328244
$ws.Range("B6").Value = 0.06657223796033994
$ws.Range("D6").Value = 0.009915014164305949
$ws.Range("F6").Value = 0.06232294617563739
$ws.Range("J6").Value = 0.2677053824362606
$ws.Range("O6").Value = 0.028328611898017
$ws.Range("Q6").Value = 0.1628895184135977
$ws.Range("R6").Value = 0.07082152974504249
$ws.Range("S6").Value = 0.3314447592067989
$ws.Range("B7").Value = 0.09731543624161074
$ws.Range("D7").Value = 0.02516778523489933
$ws.Range("F7").Value = 0.06040268456375839
$ws.Range("J7").Value = 0.151006711409396
$ws.Range("O7").Value = 0.01342281879194631
$ws.Range("Q7").Value = 0.1543624161073825
$ws.Range("R7").Value = 0.07718120805369127
$ws.Range("S7").Value = 0.4211409395973154
$ws.Range("B8").Value = 0.09254901960784313
$ws.Range("D8").Value = 0.02274509803921568
$ws.Range("E8").Value = 0.001568627450980392
$ws.Range("F8").Value = 0.0603921568627451
$ws.Range("J8").Value = 0.1184313725490196
$ws.Range("O8").Value = 0.01803921568627451
$ws.Range("Q8").Value = 0.1749019607843137
$ws.Range("R8").Value = 0.08941176470588236
$ws.Range("S8").Value = 0.4219607843137255
$ws.Range("B9").Value = 0.09122203098106713
$ws.Range("D9").Value = 0.01549053356282272
$ws.Range("F9").Value = 0.06540447504302926
$ws.Range("J9").Value = 0.1135972461273666
$ws.Range("O9").Value = 0.02581755593803787
$ws.Range("Q9").Value = 0.197934595524957
$ws.Range("R9").Value = 0.07228915662650602
$ws.Range("S9").Value = 0.4182444061962134
$ws.Range("B10").Value = 0.1093003652711436
$ws.Range("D10").Value = 0.02107333520651869
$ws.Range("E10").Value = 0.0008429334082607474
$ws.Range("F10").Value = 0.07867378477100309
$ws.Range("J10").Value = 0.1098623208766508
$ws.Range("O10").Value = 0.01376791233492554
$ws.Range("Q10").Value = 0.2020230401798258
$ws.Range("R10").Value = 0.08738409665636415
$ws.Range("S10").Value = 0.3770722112953077
$ws.Range("G11").Value = 0.1120092378752887
$ws.Range("J11").Value = 0.1027713625866051
$ws.Range("K11").Value = 0.1662817551963048
$ws.Range("L11").Value = 0.605080831408776
$ws.Range("S11").Value = 0.0138568129330254
$ws.Range("G12").Value = 0.7568555758683729
$ws.Range("J12").Value = 0.1645338208409506
$ws.Range("K12").Value = 0.01279707495429616
$ws.Range("L12").Value = 0.03290676416819013
$ws.Range("S12").Value = 0.03290676416819013
$ws.Range("F13").Value = 0.006802721088435374
$ws.Range("G13").Value = 0.673469387755102
$ws.Range("J13").Value = 0.2448979591836735
$ws.Range("S13").Value = 0.07482993197278912
$ws.Range("F15").Value = 0.01246105919003115
$ws.Range("H15").Value = 0.1526479750778816
$ws.Range("I15").Value = 0.06542056074766354
$ws.Range("J15").Value = 0.3862928348909657
$ws.Range("K15").Value = 0.06230529595015576
$ws.Range("M15").Value = 0.009345794392523364
$ws.Range("N15").Value = 0.003115264797507788
$ws.Range("O15").Value = 0.06697819314641744
$ws.Range("S15").Value = 0.2414330218068536
$ws.Range("F16").Value = 0.03225806451612903
$ws.Range("H16").Value = 0.1859582542694497
$ws.Range("I16").Value = 0.07590132827324478
$ws.Range("J16").Value = 0.3833017077798861
$ws.Range("K16").Value = 0.1252371916508539
$ws.Range("M16").Value = 0.02466793168880456
$ws.Range("N16").Value = 0.00189753320683112
$ws.Range("O16").Value = 0.05123339658444023
$ws.Range("S16").Value = 0.1195445920303605
$ws.Range("F17").Value = 0.02301587301587302
$ws.Range("H17").Value = 0.1841269841269841
$ws.Range("I17").Value = 0.1031746031746032
$ws.Range("J17").Value = 0.3714285714285714
$ws.Range("K17").Value = 0.1063492063492063
$ws.Range("M17").Value = 0.02301587301587302
$ws.Range("O17").Value = 0.0761904761904762
$ws.Range("S17").Value = 0.1126984126984127
$ws.Range("F18").Value = 0.0247787610619469
$ws.Range("H18").Value = 0.1592920353982301
$ws.Range("I18").Value = 0.09911504424778761
$ws.Range("J18").Value = 0.4194690265486726
$ws.Range("K18").Value = 0.1132743362831858
$ws.Range("M18").Value = 0.0247787610619469
$ws.Range("N18").Value = 0.001769911504424779
$ws.Range("O18").Value = 0.06902654867256637
$ws.Range("S18").Value = 0.08849557522123894
$ws.Range("F19").Value = 0.02179310344827586
$ws.Range("H19").Value = 0.2107586206896552
$ws.Range("I19").Value = 0.08744827586206896
$ws.Range("J19").Value = 0.3633103448275862
$ws.Range("K19").Value = 0.1092413793103448
$ws.Range("M19").Value = 0.02455172413793104
$ws.Range("N19").Value = 0.001931034482758621
$ws.Range("O19").Value = 0.07117241379310345
$ws.Range("S19").Value = 0.1097931034482759

Write-Host "Applied 113 cell value updates to Sheet1"
